$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp caption
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 10:22"

# Row 19: Austria -> Austria
$ws.Range("B19").Value = 13005
$ws.Range("C19").Value = 63
$ws.Range("E19").Value = 8220

# Row 24: Australia -> Australia
$ws.Range("D24").Value = 2987
$ws.Range("E24").Value = 3066
$ws.Range("F24").Value = 81

# Row 30: Chequia -> Polonia
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 5341
$ws.Range("C30").Value = 136
$ws.Range("D30").Value = 284
$ws.Range("E30").Value = 4893
$ws.Range("F30").Value = 160
$ws.Range("H30").Value = 164

# Row 31: Polonia -> Chequia
$ws.Range("A31").Value = "Chequia"
$ws.Range("B31").Value = 5335
$ws.Range("C31").Value = 23
$ws.Range("D31").Value = 243
$ws.Range("E31").Value = 4988
$ws.Range("F31").Value = 96
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 104

# Row 62: Irak -> Estonia
$ws.Range("A62").Value = "Estonia"
$ws.Range("B62").Value = 1207
$ws.Range("C62").Value = 22
$ws.Range("D62").Value = 83
$ws.Range("E62").Value = 1100
$ws.Range("F62").Value = 9
$ws.Range("H62").Value = 24

# Row 63: Estonia -> Irak
$ws.Range("A63").Value = "Irak"
$ws.Range("B63").Value = 1202
$ws.Range("D63").Value = 452
$ws.Range("E63").Value = 681
$ws.Range("F63").Value = 0
$ws.Range("H63").Value = 69

# Row 72: Bosnia y Herzegovina -> Bosnia y Herzegovina
$ws.Range("B72").Value = 841
$ws.Range("C72").Value = 37
$ws.Range("E72").Value = 711

# Row 89: Uruguay -> Oman
$ws.Range("A89").Value = "Oman"
$ws.Range("B89").Value = 457
$ws.Range("C89").Value = 38
$ws.Range("D89").Value = 109
$ws.Range("E89").Value = 346
$ws.Range("F89").Value = 3
$ws.Range("H89").Value = 2

# Row 90: Afganistan -> Uruguay
$ws.Range("A90").Value = "Uruguay"
$ws.Range("B90").Value = 456
$ws.Range("D90").Value = 192
$ws.Range("E90").Value = 257
$ws.Range("F90").Value = 14
$ws.Range("H90").Value = 7

# Row 91: Oman -> Afganistan
$ws.Range("A91").Value = "Afganistan"
$ws.Range("B91").Value = 444
$ws.Range("D91").Value = 29
$ws.Range("E91").Value = 401
$ws.Range("F91").Value = 0
$ws.Range("H91").Value = 14

# Row 162: Birmania -> Birmania
$ws.Range("D162").Value = 1
$ws.Range("E162").Value = 18

# Row 172: Namibia -> Laos
$ws.Range("A172").Value = "Laos"
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 16

# Row 173: Mongolia -> Namibia
$ws.Range("A173").Value = "Namibia"
$ws.Range("D173").Value = 3
$ws.Range("E173").Value = 13

# Row 174: Fiyi -> Mongolia
$ws.Range("A174").Value = "Mongolia"
$ws.Range("B174").Value = 16
$ws.Range("D174").Value = 4
$ws.Range("E174").Value = 12

# Row 175: Laos -> Fiyi
$ws.Range("A175").Value = "Fiyi"

# Row 183: Seychelles -> San Cristobal y Nieves
$ws.Range("A183").Value = "San Cristobal y Nieves"

# Row 184: San Cristobal y Nieves -> Seychelles
$ws.Range("A184").Value = "Seychelles"
